$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: LICETH PAOLA GOMEZ HERRERA - periodo 2206 -> 2507, valor mora 40000 -> 56940, salario 1300000 -> 1423500
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18: LICETH PAOLA GOMEZ HERRERA - periodo 2207 -> 2208, salario 1300000 -> 1423500
$ws.Range("E18").Value = "2208"
$ws.Range("G18").Value = 1423500

# Row 19: LICETH PAOLA GOMEZ HERRERA - periodo 2208 -> 2207, salario 1300000 -> 1423500
$ws.Range("E19").Value = "2207"
$ws.Range("G19").Value = 1423500

# Row 20: DIANA LUCIA ZURITA GUERRERO (2506) -> LICETH PAOLA GOMEZ HERRERA (2206), valor mora 56940 -> 40000
$ws.Range("C20").Value = "25890064"
$ws.Range("D20").Value = "LICETH PAOLA GOMEZ HERRERA"
$ws.Range("E20").Value = "2206"
$ws.Range("F20").Value = 40000

# Row 21: LICETH PAOLA GOMEZ HERRERA (2506) -> DIANA LUCIA ZURITA GUERRERO (2507), salario 1300000 -> 1423500
$ws.Range("C21").Value = "1047391352"
$ws.Range("D21").Value = "DIANA LUCIA ZURITA GUERRERO"
$ws.Range("E21").Value = "2507"
$ws.Range("G21").Value = 1423500

# Row 22: SANTIAGO QUIÑONES VERA - periodo 2506 -> 2507
$ws.Range("E22").Value = "2507"
